$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..16 (columns A..I): Hour, Lapaz, Jaro, Mandurriao, Molo,
# Diversion, Mobile SS 1, Mobile SS 2, Megaworld
$data = @(
    @(1, 20222, 5513, 11536, 16729, 0, 0, 0, 0),
    @(2, 19240, 5184, 11048, 15329, 0, 0, 0, 0),
    @(3, 18198, 4854, 10610, 14949, 0, 0, 0, 0),
    @(4, 17345, 4658, 10252, 14469, 0, 0, 0, 0),
    @(5, 17742, 4800, 10535, 14489, 0, 0, 0, 0),
    @(6, 18085, 4890, 10899, 14611, 0, 0, 0, 0),
    @(7, 17689, 4420, 10825, 14937, 0, 0, 0, 0),
    @(8, 20608, 4502, 0, 16873, 0, 0, 0, 0),
    @(9, 25888, 5559, 14619, 19334, 0, 0, 0, 0),
    @(10, 28174, 5754, 19760, 20151, 0, 0, 0, 0),
    @(11, 29407, 6128, 21229, 20872, 0, 0, 0, 0),
    @(12, 30004, 6245, 20824, 21493, 0, 0, 0, 0),
    @(13, 30232, 6142, 21177, 21767, 0, 0, 0, 0),
    @(14, 31909, 6459, 21970, 22757, 0, 0, 0, 0),
    @(15, 0, 6674, 22172, 23064, 0, 0, 0, 0)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowVals[$c]
    }
}
